# Error Calculations and Plots
# Apply the missing-data edits to Sheet1 of the active workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the two rows that are no longer present (RM 232 and SC 92) ---
# Row 26 = "RM 232", Row 28 = "SC 92" (before any deletion).
# Deleting row 26 first shifts "SC 92" up to row 27.
$ws.Rows(26).Delete()
$ws.Rows(27).Delete()

# --- Fill in / clear values in column F (and one column D) after the two-row deletion ---
$ws.Range("F2").Value = 18.03
$ws.Range("F6").Value = $null
$ws.Range("F12").Value = 17.45
$ws.Range("F14").Value = $null
$ws.Range("F20").Value = 17.73
$ws.Range("F21").Value = 16.58
$ws.Range("F22").Value = $null
$ws.Range("F23").Value = $null
$ws.Range("D30").Value = -13.6
$ws.Range("F31").Value = 17.18
$ws.Range("D32").Value = $null
$ws.Range("F33").Value = 17.53
